$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3633.9092
$ws.Range("I45").Value = 2795
$ws.Range("K45").Value = 2795
$ws.Range("M45").Value = -2418
$ws.Range("H74").Value = 4123.591
$ws.Range("J74").Value = 5576.75
$ws.Range("L74").Value = 5576.75
$ws.Range("N74").Value = -7324.75
$ws.Range("H77").Value = 4123.591
$ws.Range("J77").Value = 5576.75
$ws.Range("L77").Value = 27883.75
$ws.Range("N77").Value = -36619.75
$ws.Range("H122").Value = 1703.5834
$ws.Range("I122").Value = 1699.125
$ws.Range("J122").Value = 1712.5
$ws.Range("K122").Value = 5097.375
$ws.Range("L122").Value = 5137.5
$ws.Range("M122").Value = -2647.375
$ws.Range("N122").Value = -10037.5
$ws.Range("H132").Value = 2887.1667
$ws.Range("J132").Value = 5712.25
$ws.Range("L132").Value = 17136.75
$ws.Range("N132").Value = -22196.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4314.2
$ws.Range("I99").Value = 4022
$ws.Range("K99").Value = 4022
$ws.Range("M99").Value = -2524
$ws.Range("H107").Value = 6144.4614
$ws.Range("I107").Value = 1646.3334
$ws.Range("K107").Value = 1646.3334
$ws.Range("M107").Value = 273.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1376.1818
$ws.Range("I16").Value = 1273.625
$ws.Range("K16").Value = 1273.625
$ws.Range("M16").Value = -986.625
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 6551.04
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -705
$ws.Range("H34").Value = 6551.04
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -798
$ws.Range("H47").Value = 24999.5
$ws.Range("I47").Value = 24999.5
$ws.Range("K47").Value = 24999.5
$ws.Range("M47").Value = -24433.5
$ws.Range("H86").Value = 1540.8334
$ws.Range("I86").Value = 1449.2
$ws.Range("K86").Value = 1449.2
$ws.Range("M86").Value = -326.2
$ws.Range("H89").Value = 1540.8334
$ws.Range("I89").Value = 1449.2
$ws.Range("K89").Value = 7246
$ws.Range("M89").Value = -1630
$ws.Range("H99").Value = 3321.65
$ws.Range("I99").Value = 3084.353
$ws.Range("K99").Value = 3084.353
$ws.Range("M99").Value = -1586.353
$ws.Range("H107").Value = 352
$ws.Range("J107").Value = 268.7143
$ws.Range("L107").Value = 268.7143
$ws.Range("N107").Value = -4108.7143
$ws.Range("H113").Value = 1376.1818
$ws.Range("I113").Value = 1273.625
$ws.Range("K113").Value = 1273.625
$ws.Range("M113").Value = 896.375
$ws.Range("H122").Value = 1532
$ws.Range("J122").Value = 1478.8
$ws.Range("L122").Value = 4436.4
$ws.Range("N122").Value = -9336.4
$ws.Range("H126").Value = 3321.65
$ws.Range("I126").Value = 3084.353
$ws.Range("K126").Value = 9253.059000000001
$ws.Range("M126").Value = -6783.059000000001
$ws.Range("H132").Value = 4135.8125
$ws.Range("J132").Value = 6663
$ws.Range("L132").Value = 19989
$ws.Range("N132").Value = -25049
$ws.Range("H134").Value = 2414.3076
$ws.Range("I134").Value = 2532.1667
$ws.Range("K134").Value = 7596.500100000001
$ws.Range("M134").Value = -5061.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 110
$ws.Range("J24").Value = 110
$ws.Range("L24").Value = 330
$ws.Range("N24").Value = -790
$ws.Range("H80").Value = 4979.696
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 4922.1665
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 14766.4995
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -16638.4995
$ws.Range("H83").Value = 4979.696
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 4922.1665
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 44299.4985
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -53659.4985
$ws.Range("H109").Value = 169597.67
$ws.Range("I109").Value = 202381.6
$ws.Range("K109").Value = 607144.8
$ws.Range("M109").Value = -606104.8
$ws.Range("H114").Value = 1736.8462
$ws.Range("I114").Value = 1165
$ws.Range("J114").Value = 2227
$ws.Range("K114").Value = 3495
$ws.Range("L114").Value = 6681
$ws.Range("M114").Value = -241
$ws.Range("N114").Value = -13189
$ws.Range("H120").Value = 1761
$ws.Range("I120").Value = 875
$ws.Range("J120").Value = 3533
$ws.Range("K120").Value = 2625
$ws.Range("L120").Value = 10599
$ws.Range("M120").Value = 2213
$ws.Range("N120").Value = -20275
$ws.Range("H121").Value = 815.0909
$ws.Range("J121").Value = 1496.75
$ws.Range("L121").Value = 4490.25
$ws.Range("N121").Value = -7110.25
$ws.Range("H130").Value = 2973.2
$ws.Range("I130").Value = 2563
$ws.Range("J130").Value = 3246.6667
$ws.Range("K130").Value = 7689
$ws.Range("L130").Value = 9740.000100000001
$ws.Range("M130").Value = -2669
$ws.Range("N130").Value = -19780.0001
$ws.Range("H137").Value = 1385
$ws.Range("I137").Value = 899
$ws.Range("J137").Value = 1579.4
$ws.Range("K137").Value = 2697
$ws.Range("L137").Value = 4738.200000000001
$ws.Range("M137").Value = 2403
$ws.Range("N137").Value = -14938.2
$ws.Range("H138").Value = 6180.3335
$ws.Range("I138").Value = 2360.8333
$ws.Range("K138").Value = 7082.499899999999
$ws.Range("M138").Value = -1942.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("N5").Value = 0
$ws.Range("H80").Value = 2467.4
$ws.Range("I80").Value = 1749.5
$ws.Range("K80").Value = 1749.5
$ws.Range("M80").Value = -751.5
$ws.Range("H83").Value = 2467.4
$ws.Range("I83").Value = 1749.5
$ws.Range("K83").Value = 8747.5
$ws.Range("M83").Value = -3755.5
$ws.Range("H122").Value = 149126.17
$ws.Range("I122").Value = 209665.88
$ws.Range("K122").Value = 628997.64
$ws.Range("M122").Value = -626547.64

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H46").Value = 3649.12
$ws.Range("I46").Value = 1854.5555
$ws.Range("J46").Value = 4658.5625
$ws.Range("K46").Value = 1854.5555
$ws.Range("L46").Value = 4658.5625
$ws.Range("M46").Value = -1666.5555
$ws.Range("N46").Value = -5034.5625
$ws.Range("H61").Value = 6818.385
$ws.Range("I61").Value = 5825.5
$ws.Range("K61").Value = 5825.5
$ws.Range("M61").Value = -5623.5
$ws.Range("H93").Value = 1030.5
$ws.Range("I93").Value = 1030.5
$ws.Range("K93").Value = 1030.5
$ws.Range("M93").Value = 217.5
$ws.Range("H113").Value = 6818.385
$ws.Range("I113").Value = 5825.5
$ws.Range("K113").Value = 5825.5
$ws.Range("M113").Value = -3655.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 156730.77
$ws.Range("I2").Value = 156730.77
$ws.Range("K2").Value = 156730.77
$ws.Range("M2").Value = -156618.77
$ws.Range("H122").Value = 1385.875
$ws.Range("I122").Value = 1385.875
$ws.Range("K122").Value = 4157.625
$ws.Range("M122").Value = -1707.625
$ws.Range("H132").Value = 3268.5
$ws.Range("I132").Value = 3202
$ws.Range("K132").Value = 9606
$ws.Range("M132").Value = -7076
$ws.Range("H136").Value = 4109.727
$ws.Range("I136").Value = 2844.8667
$ws.Range("K136").Value = 8534.6001
$ws.Range("M136").Value = -5984.6001
